$wb = $excel.ActiveWorkbook

# --- "grilla de pruebas" sheet: nueva estrategia mercado quieto ---
$wsGrilla = $wb.Worksheets.Item("grilla de pruebas")

# Increase the "incremento" step value used by the projection (B6) from 15 to 30
$wsGrilla.Range("B6").Value = "30"

# Clear the manual override values in row 10 (E10/F10) so the formulas
# no longer contribute to the projection totals
$wsGrilla.Range("E10:F10").ClearContents()

# Move the selection as left by the user
$wsGrilla.Range("B5").Select()

# --- "estrategias" sheet: update the label and move the selection ---
$wsEstrategias = $wb.Worksheets.Item("estrategias")
$wsEstrategias.Range("B7").Value = "buena ganancia-mercado movido:"
$wsEstrategias.Range("B8").Select()

# estrategias is the sheet left active/selected in the workbook
$wsEstrategias.Activate()
